$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.321.98"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "3.517.90"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "596.79"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").Value = "174.40"
$ws.Range("E6").Value = "  +1.83%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("E9").Value = "  +6.44%  "

$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("D11").Value = "0.437"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "4.130.57"
$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").Value = "29.34"
$ws.Range("E14").Value = "  +4.03%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "67.279.76"
$ws.Range("E15").Value = "  +0.94%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000181"
$ws.Range("E16").Value = "  +1.95%  "

$ws.Range("D17").Value = "3.519.81"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "6.35"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").Value = "14.21"
$ws.Range("E19").Value = "  +2.30%  "

$ws.Range("D20").Value = "396.36"
$ws.Range("E20").Value = "  +2.04%  "

$ws.Range("D21").Value = "8.02"
$ws.Range("E21").Value = "  +1.00%  "

$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.540"
$ws.Range("E24").Value = "  +1.53%  "

$ws.Range("E25").Value = "  +0.85%  "

$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").Value = "6.33"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").Value = "1.47"
$ws.Range("E30").Value = "  -0.57%  "

$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("D32").Value = "23.97"
$ws.Range("E32").Value = "  +1.89%  "

$ws.Range("D33").Value = "7.40"
$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").Value = "1.68"
$ws.Range("E34").Value = "  +4.10%  "

$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").Value = "0.891"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("E37").Value = "  +1.24%  "

$ws.Range("D38").Value = "7.09"
$ws.Range("E38").Value = "  +6.77%  "

$ws.Range("D39").Value = "0.0757"
$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("D40").Value = "4.71"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("E41").Value = "  +1.67%  "

$ws.Range("D42").Value = "27.40"
$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.841.10"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.64"
$ws.Range("E44").Value = "  +4.63%  "

$ws.Range("D45").Value = "43.05"
$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("E46").Value = "  -1.33%  "

$ws.Range("D47").Value = "341.35"
$ws.Range("E47").Value = "  -4.32%  "

$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").Value = "34.37"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").Value = "6.52"
$ws.Range("E50").Value = "  +0.77%  "

$ws.Range("D51").Value = "0.851"
$ws.Range("E51").Value = "  -0.89%  "

